$p = $ppt.ActivePresentation

# --- Slide 1: Title placeholder -> "This is a title" ---
$s1 = $p.Slides.Item(1)
$titleShape1 = $s1.Shapes.Item(1)
$titleShape1.TextFrame.TextRange.Text = "This is a title"
$titleFont1 = $titleShape1.TextFrame.TextRange.Font
$titleFont1.Size = 24
$titleFont1.Name = "Calibri"
$titleFont1.Color.RGB = 0

# --- Slide 1: Body placeholder -> "test" ---
$bodyShape1 = $s1.Shapes.Item(2)
$bodyShape1.TextFrame.TextRange.Text = "test"
$bodyFont1 = $bodyShape1.TextFrame.TextRange.Font
$bodyFont1.Size = 12
$bodyFont1.Name = "Calibri"
$bodyFont1.Color.RGB = 255

# --- Slide 2: Title placeholder -> "Yes" ---
$s2 = $p.Slides.Item(2)
$titleShape2 = $s2.Shapes.Item(1)
$titleShape2.TextFrame.TextRange.Text = "Yes"
$titleFont2 = $titleShape2.TextFrame.TextRange.Font
$titleFont2.Size = 44
$titleFont2.Name = "Calibri"
$titleFont2.Color.RGB = 16711680
